$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 4 and 5 (HIT Technologies and International Parkside Products entries)
$ws.Rows("4:5").Delete()

# Update row 2 (now representing a consolidated single-company aggregate, company_name = "1")
$ws.Range("B2").Value = "'1"
$ws.Range("D2").Value = -0.0638
$ws.Range("G2").Value = 0.2150387596899225
$ws.Range("H2").Value = 0.01813953488372093
$ws.Range("I2").Value = -0.1899642401485432
$ws.Range("J2").Value = -0.1899642401485432
$ws.Range("K2").Value = -4.98
$ws.Range("L2").Value = -0.386046511627907
$ws.Range("U2").Value = 3.74
$ws.Range("V2").Value = 0.2579310344827586
$ws.Range("W2").Value = -0.343448275862069
$ws.Range("X2").Value = 0.07612369951267055
$ws.Range("Y2").Value = -0.4195719753747396
$ws.Range("Z2").Value = 0.963211769912876
$ws.Range("AA2").Value = -0.1829757919736329
$ws.Range("AB2").Value = 0.0649202514740444
$ws.Range("AC2").Value = -0.2478960434476773
$ws.Range("AD2").Value = 4.59
$ws.Range("AE2").Value = 0.07269348958103417
$ws.Range("AF2").Value = 4.662693489581034
$ws.Range("AG2").Value = 0.9226934895810341
$ws.Range("AH2").Value = 0.2433214042752493
$ws.Range("AI2").Value = 0.3280654362242552
$ws.Range("AJ2").Value = 0.05982700040069976
$ws.Range("AK2").Value = 0.08810469727763864
$ws.Range("AL2").Value = 0.205
$ws.Range("AM2").Value = 0.174
$ws.Range("AN2").Value = -9.849785407725321
$ws.Range("AO2").Value = -12.09756097560976
$ws.Range("AP2").Value = -1.980028947598786
$ws.Range("AQ2").Value = -14.25287356321839

# Update row 3 (company renamed to D-BOX Technologies Inc., values now mirror row 2)
$ws.Range("B3").Value = "D-BOX Technologies Inc. (TSX:DBO)"
$ws.Range("D3").Value = -0.0638
$ws.Range("G3").Value = 0.2150387596899225
$ws.Range("H3").Value = 0.01813953488372093
$ws.Range("I3").Value = -0.1899642401485432
$ws.Range("J3").Value = -0.1899642401485432
$ws.Range("K3").Value = -4.98
$ws.Range("L3").Value = -0.386046511627907
$ws.Range("U3").Value = 3.74
$ws.Range("V3").Value = 0.2579310344827586
$ws.Range("W3").Value = -0.343448275862069
$ws.Range("X3").Value = 0.07612369951267055
$ws.Range("Y3").Value = -0.4195719753747396
$ws.Range("Z3").Value = 0.963211769912876
$ws.Range("AA3").Value = -0.1829757919736329
$ws.Range("AB3").Value = 0.0649202514740444
$ws.Range("AC3").Value = -0.2478960434476773
$ws.Range("AD3").Value = 4.59
$ws.Range("AE3").Value = 0.07269348958103417
$ws.Range("AF3").Value = 4.662693489581034
$ws.Range("AG3").Value = 0.9226934895810341
$ws.Range("AH3").Value = 0.2433214042752493
$ws.Range("AI3").Value = 0.3280654362242552
$ws.Range("AJ3").Value = 0.05982700040069976
$ws.Range("AK3").Value = 0.08810469727763864
$ws.Range("AL3").Value = 0.205
$ws.Range("AM3").Value = 0.174
$ws.Range("AN3").Value = -9.849785407725321
$ws.Range("AO3").Value = -12.09756097560976
$ws.Range("AP3").Value = -1.980028947598786
$ws.Range("AQ3").Value = -14.25287356321839
